$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Concepts to Learn PE")
$ws.Select()

# Write new cells in the exact order the strings first appear in the
# target workbook so that shared-string indices line up (349..361).
$ws.Range("N4").Value  = "AWR report,Thread Dump analsyis ,Heap Dump analysis"
$ws.Range("N5").Value  = "Performnce BottleNecks"
$ws.Range("N6").Value  = "App Dynamics"
$ws.Range("N7").Value  = "SAP protocol,Ajax truclient"
$ws.Range("N9").Value  = "Adobe analytics,Splunk"
$ws.Range("N11").Value = "Jenkins"
$ws.Range("N12").Value = "Master slave architeture/jmeter"
$ws.Range("N13").Value = "Prometheus ,grafana,infllux"
$ws.Range("N8").Value  = "NeoLoad, Load Runner"
$ws.Range("N10").Value = "K8S and Docker"
$ws.Range("N14").Value = "git & GitHub"

$ws.Range("N3").Value = "Top Items to Focus"
$ws.Range("N3").Font.Bold = $true

$ws.Range("J1").Value = "For 5 years Exp"
$ws.Range("J1").Font.Bold = $true

$ws.Range("Q16").Select()
$excel.ActiveWindow.Zoom = 80
